# Auto-generated update of H:N leve-profit columns across job sheets
# (values refreshed by the scheduled market-price runner)
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1485.1333
$ws.Range("I40").Value = 1616.091
$ws.Range("K40").Value = 1616.091
$ws.Range("M40").Value = -1441.091
$ws.Range("H58").Value = 625630.4
$ws.Range("I58").Value = 186.55556
$ws.Range("J58").Value = 1429772.4
$ws.Range("K58").Value = 559.66668
$ws.Range("L58").Value = 4289317.199999999
$ws.Range("M58").Value = -409.66668
$ws.Range("N58").Value = -4289617.199999999
$ws.Range("H111").Value = 55558290
$ws.Range("I111").Value = 2184.2307
$ws.Range("J111").Value = 200004160
$ws.Range("K111").Value = 6552.6921
$ws.Range("L111").Value = 600012480
$ws.Range("M111").Value = -3485.6921
$ws.Range("N111").Value = -600018614
$ws.Range("H128").Value = 15300.869
$ws.Range("J128").Value = 15300.869
$ws.Range("L128").Value = 15300.869
$ws.Range("N128").Value = -25260.869
$ws.Range("H129").Value = 1034.7192
$ws.Range("I129").Value = 485.5
$ws.Range("J129").Value = 1099.3334
$ws.Range("K129").Value = 1456.5
$ws.Range("L129").Value = 3298.0002
$ws.Range("M129").Value = 3543.5
$ws.Range("N129").Value = -13298.0002
$ws.Range("H135").Value = 862.3226
$ws.Range("I135").Value = 654.0741
$ws.Range("J135").Value = 2268
$ws.Range("K135").Value = 5886.6669
$ws.Range("L135").Value = 20412
$ws.Range("M135").Value = -3351.6669
$ws.Range("N135").Value = -25482
$ws.Range("H137").Value = 1493.1794
$ws.Range("I137").Value = 1521.6666
$ws.Range("J137").Value = 1459.9445
$ws.Range("K137").Value = 4564.9998
$ws.Range("L137").Value = 4379.833500000001
$ws.Range("M137").Value = -2014.9998
$ws.Range("N137").Value = -9479.833500000001
$ws.Range("H138").Value = 4313.983
$ws.Range("J138").Value = 6008.294
$ws.Range("L138").Value = 18024.882
$ws.Range("N138").Value = -28304.882
$ws.Range("H141").Value = 2480
$ws.Range("I141").Value = 2560
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 7680
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -2500
$ws.Range("N141").Value = -16360

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5157.3394
$ws.Range("I32").Value = 3463.3696
$ws.Range("K32").Value = 3463.3696
$ws.Range("M32").Value = -3176.3696
$ws.Range("H110").Value = 67806.55499999999
$ws.Range("I110").Value = 120506.6
$ws.Range("J110").Value = 1931.5
$ws.Range("K110").Value = 120506.6
$ws.Range("L110").Value = 1931.5
$ws.Range("M110").Value = -118461.6
$ws.Range("N110").Value = -6021.5
$ws.Range("H122").Value = 17546582
$ws.Range("I122").Value = 37038916
$ws.Range("J122").Value = 3481.2
$ws.Range("K122").Value = 111116748
$ws.Range("L122").Value = 10443.6
$ws.Range("M122").Value = -111114298
$ws.Range("N122").Value = -15343.6
$ws.Range("H132").Value = 2844.5417
$ws.Range("I132").Value = 1475.4445
$ws.Range("K132").Value = 4426.333500000001
$ws.Range("M132").Value = -1896.333500000001

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1417.6757
$ws.Range("I99").Value = 888.6667
$ws.Range("J99").Value = 1778.3636
$ws.Range("K99").Value = 888.6667
$ws.Range("L99").Value = 1778.3636
$ws.Range("M99").Value = 609.3333
$ws.Range("N99").Value = -4774.3636
$ws.Range("H106").Value = 25000
$ws.Range("J106").Value = 25000
$ws.Range("L106").Value = 25000
$ws.Range("N106").Value = -27524
$ws.Range("H107").Value = 1102.375
$ws.Range("I107").Value = 1136.6666
$ws.Range("J107").Value = 999.5
$ws.Range("K107").Value = 1136.6666
$ws.Range("L107").Value = 999.5
$ws.Range("M107").Value = 783.3334
$ws.Range("N107").Value = -4839.5
$ws.Range("H134").Value = 2901.1333
$ws.Range("I134").Value = 1578.1111
$ws.Range("J134").Value = 4885.6665
$ws.Range("K134").Value = 4734.3333
$ws.Range("L134").Value = 14656.9995
$ws.Range("M134").Value = -2199.3333
$ws.Range("N134").Value = -19726.9995

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 31303
$ws.Range("I23").Value = 14899
$ws.Range("J23").Value = 39505
$ws.Range("K23").Value = 14899
$ws.Range("L23").Value = 39505
$ws.Range("M23").Value = -14659
$ws.Range("N23").Value = -39985
$ws.Range("H27").Value = 31303
$ws.Range("I27").Value = 14899
$ws.Range("J27").Value = 39505
$ws.Range("K27").Value = 14899
$ws.Range("L27").Value = 39505
$ws.Range("M27").Value = -14707
$ws.Range("N27").Value = -39889
$ws.Range("H92").Value = 30197.5
$ws.Range("J92").Value = 30197.5
$ws.Range("L92").Value = 30197.5
$ws.Range("N92").Value = -35189.5

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 7507.524
$ws.Range("J39").Value = 2882.75
$ws.Range("L39").Value = 8648.25
$ws.Range("N39").Value = -9236.25
$ws.Range("H55").Value = 4900
$ws.Range("J55").Value = 4900
$ws.Range("L55").Value = 14700
$ws.Range("N55").Value = -15054
$ws.Range("H68").Value = 491383.16
$ws.Range("I68").Value = 1488088.4
$ws.Range("J68").Value = 1773.5614
$ws.Range("K68").Value = 4464265.199999999
$ws.Range("L68").Value = 5320.6842
$ws.Range("M68").Value = -4463454.199999999
$ws.Range("N68").Value = -6942.6842
$ws.Range("H71").Value = 491383.16
$ws.Range("I71").Value = 1488088.4
$ws.Range("J71").Value = 1773.5614
$ws.Range("K71").Value = 13392795.6
$ws.Range("L71").Value = 15962.0526
$ws.Range("M71").Value = -13388739.6
$ws.Range("N71").Value = -24074.0526
$ws.Range("H131").Value = 4423.606
$ws.Range("J131").Value = 4677.387
$ws.Range("L131").Value = 14032.161
$ws.Range("N131").Value = -24112.161
$ws.Range("H141").Value = 4280
$ws.Range("I141").Value = 826.6667
$ws.Range("J141").Value = 25000
$ws.Range("K141").Value = 2480.0001
$ws.Range("L141").Value = 75000
$ws.Range("M141").Value = 2699.9999
$ws.Range("N141").Value = -85360

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17693.572
$ws.Range("I70").Value = 19975.834
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 19975.834
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -19705.834
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 17693.572
$ws.Range("I73").Value = 19975.834
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 19975.834
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -19039.834
$ws.Range("N73").Value = -5872
$ws.Range("H80").Value = 1973
$ws.Range("I80").Value = 1977.5
$ws.Range("J80").Value = 1965.5
$ws.Range("K80").Value = 1977.5
$ws.Range("L80").Value = 1965.5
$ws.Range("M80").Value = -979.5
$ws.Range("N80").Value = -3961.5
$ws.Range("H83").Value = 1973
$ws.Range("I83").Value = 1977.5
$ws.Range("J83").Value = 1965.5
$ws.Range("K83").Value = 9887.5
$ws.Range("L83").Value = 9827.5
$ws.Range("M83").Value = -4895.5
$ws.Range("N83").Value = -19811.5
$ws.Range("H102").Value = 903.53845
$ws.Range("I102").Value = 812.1667
$ws.Range("K102").Value = 812.1667
$ws.Range("M102").Value = 809.8333
$ws.Range("H105").Value = 50335.5
$ws.Range("J105").Value = 50335.5
$ws.Range("L105").Value = 50335.5
$ws.Range("N105").Value = -57323.5
$ws.Range("H132").Value = 3843.9443
$ws.Range("I132").Value = 3396
$ws.Range("J132").Value = 4016.2307
$ws.Range("K132").Value = 10188
$ws.Range("L132").Value = 12048.6921
$ws.Range("M132").Value = -7658
$ws.Range("N132").Value = -17108.6921

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3599.111
$ws.Range("I16").Value = 2878
$ws.Range("K16").Value = 2878
$ws.Range("M16").Value = -2708
$ws.Range("H22").Value = 886.5789
$ws.Range("I22").Value = 760.1
$ws.Range("K22").Value = 760.1
$ws.Range("M22").Value = -465.1
$ws.Range("H27").Value = 886.5789
$ws.Range("I27").Value = 760.1
$ws.Range("K27").Value = 760.1
$ws.Range("M27").Value = -653.1
$ws.Range("H55").Value = 363.16666
$ws.Range("I55").Value = 363.16666
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 363.16666
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -190.16666
$ws.Range("N55").ClearContents()
$ws.Range("H122").Value = 7939614
$ws.Range("I122").Value = 13890499
$ws.Range("J122").Value = 5101
$ws.Range("K122").Value = 41671497
$ws.Range("L122").Value = 15303
$ws.Range("M122").Value = -41669047
$ws.Range("N122").Value = -20203
$ws.Range("H132").Value = 4243.316
$ws.Range("I132").Value = 5134.074
$ws.Range("J132").Value = 3441.6333
$ws.Range("K132").Value = 15402.222
$ws.Range("L132").Value = 10324.8999
$ws.Range("M132").Value = -12872.222
$ws.Range("N132").Value = -15384.8999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("M56").ClearContents()
$ws.Range("N56").ClearContents()
$ws.Range("H81").Value = 1636.2727
$ws.Range("I81").Value = 1999.8572
$ws.Range("K81").Value = 3999.7144
$ws.Range("M81").Value = -2938.7144
$ws.Range("H84").Value = 1636.2727
$ws.Range("I84").Value = 1999.8572
$ws.Range("K84").Value = 19998.572
$ws.Range("M84").Value = -14694.572
$ws.Range("H92").Value = 25550
$ws.Range("J92").Value = 25550
$ws.Range("L92").Value = 25550
$ws.Range("N92").Value = -30542
$ws.Range("H107").Value = 533.3570999999999
$ws.Range("I107").Value = 531.7143
$ws.Range("J107").Value = 535
$ws.Range("K107").Value = 1595.1429
$ws.Range("L107").Value = 1605
$ws.Range("M107").Value = 324.8571000000002
$ws.Range("N107").Value = -5445
$ws.Range("H122").Value = 31183.97
$ws.Range("I122").Value = 38398
$ws.Range("J122").Value = 3358.4285
$ws.Range("K122").Value = 115194
$ws.Range("L122").Value = 10075.2855
$ws.Range("M122").Value = -112744
$ws.Range("N122").Value = -14975.2855
$ws.Range("H132").Value = 1756.7358
$ws.Range("I132").Value = 1302.7333
$ws.Range("J132").Value = 2348.913
$ws.Range("K132").Value = 3908.199900000001
$ws.Range("L132").Value = 7046.739
$ws.Range("M132").Value = -1378.199900000001
$ws.Range("N132").Value = -12106.739

